$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $CellRef, $Text) {
    $r = $Sheet.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Text
    $r.Style = "Normal"
}

Set-TextValue $ws 'D2' '286.11'
Set-TextValue $ws 'G2' '2'
Set-TextValue $ws 'D3' '21.06'
Set-TextValue $ws 'G3' '2'
Set-TextValue $ws 'D4' '6.436'
Set-TextValue $ws 'G4' '2'
Set-TextValue $ws 'D5' '0.06226'
Set-TextValue $ws 'G5' '2'
Set-TextValue $ws 'D6' '3.588'
Set-TextValue $ws 'G6' '2'
Set-TextValue $ws 'B7' 'FTXToken'
Set-TextValue $ws 'C7' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws 'D7' '1.540'
Set-TextValue $ws 'E7' '6FTXTokenFTT'
Set-TextValue $ws 'G7' '2'
Set-TextValue $ws 'B8' 'KuCoinToken'
Set-TextValue $ws 'C8' 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue $ws 'D8' '6.574'
Set-TextValue $ws 'E8' '7KuCoinTokenKCS'
Set-TextValue $ws 'G8' '2'
Set-TextValue $ws 'D9' '0.8230'
Set-TextValue $ws 'G9' '2'
Set-TextValue $ws 'D10' '0.01407'
Set-TextValue $ws 'G10' '2'
Set-TextValue $ws 'D11' '0.1660'
Set-TextValue $ws 'G11' '2'
Set-TextValue $ws 'D12' '0.08548'
Set-TextValue $ws 'G12' '2'
Set-TextValue $ws 'D13' '0.03530'
Set-TextValue $ws 'G13' '2'
Set-TextValue $ws 'D14' '0.03221'
Set-TextValue $ws 'G14' '2'
Set-TextValue $ws 'D15' '0.09191'
Set-TextValue $ws 'G15' '2'
Set-TextValue $ws 'G16' '2'
Set-TextValue $ws 'D17' '0.001651'
Set-TextValue $ws 'G17' '2'
Set-TextValue $ws 'D18' '0.04755'
Set-TextValue $ws 'G18' '2'
Set-TextValue $ws 'D19' '0.006265'
Set-TextValue $ws 'G19' '2'
Set-TextValue $ws 'D20' '0.006231'
Set-TextValue $ws 'G20' '2'
Set-TextValue $ws 'G21' '2'
Set-TextValue $ws 'D22' '0.0001602'
Set-TextValue $ws 'G22' '2'
Set-TextValue $ws 'D23' '3.824'
Set-TextValue $ws 'G23' '2'
Set-TextValue $ws 'D24' '2.334'
Set-TextValue $ws 'G24' '2'
Set-TextValue $ws 'G25' '2'
Set-TextValue $ws 'D26' '0.1224'
Set-TextValue $ws 'G26' '2'
Set-TextValue $ws 'G27' '2'
Set-TextValue $ws 'G28' '2'
Set-TextValue $ws 'G29' '2'
Set-TextValue $ws 'G30' '2'
Set-TextValue $ws 'G31' '2'
Set-TextValue $ws 'G32' '2'
Set-TextValue $ws 'G33' '2'
Set-TextValue $ws 'G34' '2'
Set-TextValue $ws 'G35' '2'
Set-TextValue $ws 'G36' '2'
Set-TextValue $ws 'G37' '2'
Set-TextValue $ws 'G38' '2'
Set-TextValue $ws 'G39' '2'
Set-TextValue $ws 'D40' '0.04738'
Set-TextValue $ws 'G40' '2'
Set-TextValue $ws 'B41' 'KickToken'
Set-TextValue $ws 'C41' 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws 'D41' '0.007149'
Set-TextValue $ws 'E41' '40KickTokenKICK'
Set-TextValue $ws 'G41' '2'
Set-TextValue $ws 'D42' '0.004506'
Set-TextValue $ws 'G42' '2'
Set-TextValue $ws 'B43' 'BKEXToken'
Set-TextValue $ws 'C43' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws 'D43' '0.1108'
Set-TextValue $ws 'E43' '42BKEXTokenBKK'
Set-TextValue $ws 'G43' '2'
Set-TextValue $ws 'D44' '0.01170'
Set-TextValue $ws 'G44' '2'
Set-TextValue $ws 'D45' '0.00006933'
Set-TextValue $ws 'G45' '2'
Set-TextValue $ws 'D46' '0.00000000751'
Set-TextValue $ws 'G46' '2'
Set-TextValue $ws 'D47' '0.9022'
Set-TextValue $ws 'G47' '2'
Set-TextValue $ws 'D48' '0.002892'
Set-TextValue $ws 'G48' '2'
Set-TextValue $ws 'D49' '0.00001402'
Set-TextValue $ws 'G49' '2'
Set-TextValue $ws 'D50' '0.01242'
Set-TextValue $ws 'G50' '2'
Set-TextValue $ws 'G51' '2'
